$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: merge Effort + Additional Effort into a single Effort value (3.75 + 0.25 -> 4)
$ws.Range("B14").Value = 4
$ws.Range("C14").ClearContents()

# Row 37: merge Effort + Additional Effort into a single Effort value (1.75 + 0.25 -> 2)
$ws.Range("B37").Value = 2
$ws.Range("C37").ClearContents()

# New row 41: new entry for Test case tc08
$ws.Range("A41").Value = 41227
$ws.Range("B41").Value = 2.5
$ws.Range("D41").Value = "Test case tc08 put to operation but still shows some bad behavior"

# Update selection / view to match the saved state of the workbook
$ws.Range("B22").Select()
